$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$emails = @(
  "newmaker.plus0@gmail.com",
  "1050223@gm.yhsh.tn.edu.tw",
  "810002@gm.yhsh.tn.edu.tw",
  "C14116132@gs.ncku.edu.tw",
  "c14116132@gmail.com",
  "mr.871945@gmail.com"
)

for ($i = 0; $i -lt $emails.Length; $i++) {
  $row = $i + 1
  $addrCell = $ws.Cells.Item($row, 1)
  [void]$ws.Hyperlinks.Add($addrCell, "mailto:" + $emails[$i])
  $addrCell.Value = $emails[$i]
  $addrCell.Style = "超連結"
  $ws.Cells.Item($row, 2).Value = "asd1016101610"
}

[void]$ws.Range("B7").Select()
